# Update EUR->ARS rate: 2025-10-10T15:21:25Z
# Appends a new row to the rate-history sheet with the latest quote.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 69
$rangeAddr = "A" + $newRow + ":B" + $newRow

# Columns A/B hold date-like / time-like text ("2025-10-10", "15:21:25").
# Plain assignment would let Excel's smart-entry coerce these into
# date/time serial numbers, so force text interpretation first, then
# restore the default "Normal" style afterwards so the new cells end up
# unformatted, matching every other row already on the sheet.
$ws.Range($rangeAddr).NumberFormat = "@"

$ws.Cells.Item($newRow, 1).Value = "2025-10-10"
$ws.Cells.Item($newRow, 2).Value = "15:21:25"
$ws.Cells.Item($newRow, 3).Value = "1.00 EUR = 1,750.2781"

$ws.Range($rangeAddr).Style = "Normal"
